# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before the
#    "总计" summary sheet) and fill it with the per-fund holding detail for
#    Q1-2022.
# 2. Update the "总计" (totals) summary sheet: add a new top row for
#    "2022-Q1" (14 funds, 2.42 billion RMB held) and push the existing
#    "2021-Q4" row down by one, renumbering its index.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q1" worksheet, positioned after "2021-Q4"
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Re-fetch everything by stable name/index from here on -- this runtime
# resolves previously-captured worksheet variables by their *position*,
# which shifts once sheets are inserted, so we look sheets up again fresh
# whenever the sheet collection may have changed.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("2022-Q1")

# Copy the header-row formatting (bold / border / centered) from the
# "2021-Q4" sheet onto the matching header cells of the new sheet.
$q4.Range("B1:H1").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Copy the row-index column formatting (column A) from "2021-Q4" as well.
$q4.Range("A2").Copy() | Out-Null
$q1.Range("A2:A15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Columns B-G hold text (fund codes / names / numbers-as-text, matching the
# "2021-Q4" sheet's own representation), so force a text number format
# before writing the values to avoid Excel auto-converting them to numbers
# (which would, e.g., strip leading zeros from fund codes).
$q1.Range("B2:G15").NumberFormat = "@"

$rows = @(
    @(0,  "009646", "南方核心成长混合A",               "19.33", "84.70", "3.57", "0.6901", 7),
    @(1,  "202011", "南方优选价值混合A",               "11.61", "84.08", "3.58", "0.4156", 8),
    @(2,  "010132", "南方创新成长混合A",               "11.16", "85.17", "3.62", "0.4040", 7),
    @(3,  "009681", "南方创新精选一年定期开放混合A",   "11.24", "92.15", "3.52", "0.3956", 8),
    @(4,  "009647", "南方核心成长混合C",               "3.19",  "84.70", "3.57", "0.1139", 7),
    @(5,  "009682", "南方创新精选一年定期开放混合C",   "3.09",  "92.15", "3.52", "0.1088", 8),
    @(6,  "005729", "南方人工智能主题混合",             "2.12",  "83.28", "4.89", "0.1037", 4),
    @(7,  "002577", "南方新兴龙头灵活配置混合",         "2.09",  "78.47", "4.79", "0.1001", 2),
    @(8,  "010133", "南方创新成长混合C",               "2.06",  "85.17", "3.62", "0.0746", 7),
    @(9,  "013903", "国泰君安信息行业混合",             "0.25",  "84.06", "2.88", "0.0072", 9),
    @(10, "006539", "南方优选价值混合C",               "0.12",  "84.08", "3.58", "0.0043", 8),
    @(11, "960020", "南方优选价值混合H",               "0.08",  "84.08", "3.58", "0.0029", 8),
    @(12, "004976", "华润元大景泰混合A",               "1.79",  "37.61", "0.11", "0.0020", 10),
    @(13, "004977", "华润元大景泰混合C",               "1.79",  "37.61", "0.11", "0.0020", 10)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet -- push the existing 2021-Q4 row
# down to row 3 (index 1) and write the new 2022-Q1 totals into row 2
# (index 0).
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2021-Q4"
$zj.Cells.Item(3, 3).Value = 2
$zj.Cells.Item(3, 4).Value = 0.22

# copy the row-index style (bold/border/centered) from A2 down onto the
# newly-used A3 cell so both index cells look the same.
$zj.Range("A2").Copy() | Out-Null
$zj.Range("A3").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 14
$zj.Cells.Item(2, 4).Value = 2.42

Write-Output "2022-Q1 sheet and 总计 summary updated"
